$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 'fragen'
$ws.Range("B3").Value = 'loben'
$ws.Range("B4").Value = 'brauen'
$ws.Range("B5").Value = 'schwören'
$ws.Range("B6").Value = 'arten'
$ws.Range("B7").Value = 'bitten'
$ws.Range("B8").Value = 'reiten'
$ws.Range("B9").Value = 'bellen'
$ws.Range("B10").Value = 'mauern'
$ws.Range("B11").Value = 'grüßen'
$ws.Range("B12").Value = 'sprechen'
$ws.Range("B13").Value = 'stehlen'
$ws.Range("B14").Value = 'ärgern'
$ws.Range("B15").Value = 'herrschen'
$ws.Range("B16").Value = 'quälen'
$ws.Range("B17").Value = 'saufen'
$ws.Range("B18").Value = 'filmen'
$ws.Range("B19").Value = 'dringen'
$ws.Range("B20").Value = 'sichern'
$ws.Range("B21").Value = 'baden'
$ws.Range("B22").Value = 'wüten'
$ws.Range("B23").Value = 'klettern'
$ws.Range("B24").Value = 'küssen'
$ws.Range("B25").Value = 'scheinen'
$ws.Range("B26").Value = 'steuern'
$ws.Range("B27").Value = 'äußern'
$ws.Range("B28").Value = 'enden'
$ws.Range("B29").Value = 'spinnen'
$ws.Range("B30").Value = 'retten'
$ws.Range("B31").Value = 'hören'
$ws.Range("B32").Value = 'sinken'
$ws.Range("B33").Value = 'lieben'
$ws.Range("B34").Value = 'fischen'
$ws.Range("B35").Value = 'rufen'
$ws.Range("B36").Value = 'wärmen'
$ws.Range("B37").Value = 'führen'
$ws.Range("B38").Value = 'lehnen'
$ws.Range("B39").Value = 'schreiten'
$ws.Range("B40").Value = 'zünden'
$ws.Range("B41").Value = 'sitzen'
$ws.Range("B42").Value = 'ahnen'
$ws.Range("B43").Value = 'flüchten'
$ws.Range("B44").Value = 'fahren'
$ws.Range("B45").Value = 'sprengen'
$ws.Range("B46").Value = 'schnellen'
$ws.Range("B47").Value = 'wenden'
$ws.Range("B48").Value = 'lesen'
$ws.Range("B49").Value = 'kürzen'
$ws.Range("B50").Value = 'beten'
$ws.Range("B51").Value = 'öffnen'
$ws.Range("B52").Value = 'jubeln'
$ws.Range("B53").Value = 'sterben'
$ws.Range("B54").Value = 'streichen'
$ws.Range("B55").Value = 'planen'
$ws.Range("B56").Value = 'weigern'
$ws.Range("B57").Value = 'biegen'
$ws.Range("B58").Value = 'machen'
$ws.Range("B59").Value = 'leeren'
$ws.Range("B60").Value = 'stecken'
$ws.Range("B61").Value = 'reizen'
$ws.Range("B62").Value = 'trennen'
$ws.Range("B63").Value = 'bilden'
$ws.Range("B64").Value = 'deuten'
$ws.Range("B65").Value = 'wachsen'
$ws.Range("B66").Value = 'drehen'
$ws.Range("B67").Value = 'schätzen'
$ws.Range("B68").Value = 'kehren'
$ws.Range("B69").Value = 'suchen'
$ws.Range("B70").Value = 'albern'
$ws.Range("B71").Value = 'folgen'
$ws.Range("B72").Value = 'hoffen'
$ws.Range("B73").Value = 'zeigen'
$ws.Range("B74").Value = 'stürmen'
$ws.Range("B75").Value = 'altern'
$ws.Range("B76").Value = 'kümmern'
$ws.Range("B77").Value = 'streifen'
$ws.Range("B78").Value = 'sparen'
$ws.Range("B79").Value = 'mögen'
$ws.Range("B80").Value = 'schrecken'
$ws.Range("B81").Value = 'segeln'
$ws.Range("B82").Value = 'hauen'
$ws.Range("B83").Value = 'scheitern'
$ws.Range("B84").Value = 'schmecken'
$ws.Range("B85").Value = 'kranken'
$ws.Range("B86").Value = 'knarren'
$ws.Range("B87").Value = 'freuen'
$ws.Range("B88").Value = 'graben'
$ws.Range("B89").Value = 'bremsen'
$ws.Range("B90").Value = 'warnen'
$ws.Range("B91").Value = 'irren'
$ws.Range("B92").Value = 'kosten'
$ws.Range("B93").Value = 'achten'
$ws.Range("B94").Value = 'nennen'
$ws.Range("B95").Value = 'geben'
$ws.Range("B96").Value = 'bluten'
$ws.Range("B97").Value = 'schlucken'
$ws.Range("B98").Value = 'bauen'
$ws.Range("B99").Value = 'kriegen'
$ws.Range("B100").Value = 'dienen'
$ws.Range("B101").Value = 'töten'
$ws.Range("B102").Value = 'werfen'
$ws.Range("B103").Value = 'bergen'
$ws.Range("B104").Value = 'wehtun'
$ws.Range("B105").Value = 'wandern'
$ws.Range("B106").Value = 'pfeifen'
$ws.Range("B107").Value = 'haben'
$ws.Range("B108").Value = 'ehren'
$ws.Range("B109").Value = 'strahlen'
$ws.Range("B110").Value = 'backen'
$ws.Range("B111").Value = 'testen'
$ws.Range("B112").Value = 'tollen'
$ws.Range("B113").Value = 'fangen'
$ws.Range("B114").Value = 'erben'
$ws.Range("B115").Value = 'siegen'
$ws.Range("B116").Value = 'schlagen'
$ws.Range("B117").Value = 'heben'
$ws.Range("B118").Value = 'sammeln'
$ws.Range("B119").Value = 'spüren'
$ws.Range("B120").Value = 'sorgen'
$ws.Range("B121").Value = 'lösen'
$ws.Range("B122").Value = 'tropfen'
$ws.Range("B123").Value = 'helfen'
$ws.Range("B124").Value = 'tanzen'
$ws.Range("B125").Value = 'lockern'
$ws.Range("B126").Value = 'tauschen'
$ws.Range("B127").Value = 'sperren'
$ws.Range("B128").Value = 'trauen'
$ws.Range("B129").Value = 'regeln'
$ws.Range("B130").Value = 'kichern'
$ws.Range("B131").Value = 'knien'
$ws.Range("B132").Value = 'trösten'
$ws.Range("B133").Value = 'jagen'
$ws.Range("B134").Value = 'fließen'
$ws.Range("B135").Value = 'posten'
$ws.Range("B136").Value = 'weichen'
$ws.Range("B137").Value = 'stammen'
$ws.Range("B138").Value = 'hassen'
$ws.Range("B139").Value = 'starten'
$ws.Range("B140").Value = 'betteln'
$ws.Range("B141").Value = 'pflanzen'
$ws.Range("B142").Value = 'feiern'
$ws.Range("B143").Value = 'wirken'
$ws.Range("B144").Value = 'schleppen'
$ws.Range("B145").Value = 'ziehen'
$ws.Range("B146").Value = 'liegen'
$ws.Range("B147").Value = 'seufzen'
$ws.Range("B148").Value = 'lohnen'
$ws.Range("B149").Value = 'werden'
$ws.Range("B150").Value = 'zögern'
$ws.Range("B151").Value = 'decken'
$ws.Range("B152").Value = 'holen'
$ws.Range("B153").Value = 'zielen'
$ws.Range("B154").Value = 'nähen'
$ws.Range("B155").Value = 'schenken'
$ws.Range("B156").Value = 'wundern'
$ws.Range("B157").Value = 'münzen'
$ws.Range("B158").Value = 'treiben'
$ws.Range("B159").Value = 'rühren'
$ws.Range("B160").Value = 'schwingen'
$ws.Range("B161").Value = 'stillen'
$ws.Range("B162").Value = 'gelten'
$ws.Range("B163").Value = 'brauchen'
$ws.Range("B164").Value = 'platzen'
$ws.Range("B165").Value = 'liefern'
$ws.Range("B166").Value = 'spielen'
$ws.Range("B167").Value = 'schwächen'
$ws.Range("B168").Value = 'weinen'
$ws.Range("B169").Value = 'hupen'
$ws.Range("B170").Value = 'lügen'
$ws.Range("B171").Value = 'klingen'
$ws.Range("B172").Value = 'runden'
$ws.Range("B173").Value = 'schulden'
$ws.Range("B174").Value = 'malen'
$ws.Range("B175").Value = 'leugnen'
$ws.Range("B176").Value = 'greifen'
$ws.Range("B177").Value = 'flehen'
$ws.Range("B178").Value = 'schultern'
$ws.Range("B179").Value = 'fällen'
$ws.Range("B180").Value = 'ändern'
$ws.Range("B181").Value = 'heilen'
$ws.Range("B182").Value = 'tragen'
$ws.Range("B183").Value = 'räumen'
$ws.Range("B184").Value = 'schauen'
$ws.Range("B185").Value = 'tauchen'
$ws.Range("B186").Value = 'fallen'
$ws.Range("B187").Value = 'heulen'
$ws.Range("B188").Value = 'plaudern'
$ws.Range("B189").Value = 'gründen'
$ws.Range("B190").Value = 'rasen'
$ws.Range("B191").Value = 'formen'
$ws.Range("B192").Value = 'boxen'
$ws.Range("B193").Value = 'foltern'
